# Applies the cryptos.xlsx data-refresh described in the commit
# "Updated cryptos list on Wed Apr  5 14:49:40 UTC 2023 with GitHub Actions".
# Rewrites the Price (D) and Volume(1h) (E) columns with refreshed scrape
# values, and swaps two row pairs (Solana/WrappedEther, Aptos/ICP) back to
# their updated rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.240.44"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3
$ws.Range("D3").Value = "1.905.82"
$ws.Range("E3").Value = "  +1.76%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'314.38"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7
$ws.Range("D7").Value = "'0.5084"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8
$ws.Range("D8").Value = "'0.3929"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.09641"
$ws.Range("E9").Value = "  -2.21%  "

# Row 10
$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  -0.08%  "

# Row 11
$ws.Range("D11").Value = "'42.00"
$ws.Range("E11").Value = "  +1.80%  "

# Row 12
$ws.Range("D12").Value = "'6.430"
$ws.Range("E12").Value = "  -0.74%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.96"
$ws.Range("E13").Value = "  -0.22%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.921.21"
$ws.Range("E14").Value = "  +2.66%  "

# Row 15
$ws.Range("E15").Value = "  -0.95%  "

# Row 16
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("E17").Value = "  -0.93%  "

# Row 18
$ws.Range("D18").Value = "'92.64"
$ws.Range("E18").Value = "  -1.02%  "

# Row 19
$ws.Range("D19").Value = "'0.06638"
$ws.Range("E19").Value = "  +0.03%  "

# Row 20
$ws.Range("D20").Value = "'17.99"
$ws.Range("E20").Value = "  +3.20%  "

# Row 21
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").Value = "'6.226"
$ws.Range("E22").Value = "  +1.64%  "

# Row 23
$ws.Range("D23").Value = "28.295.01"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").Value = "'11.31"

# Row 25
$ws.Range("D25").Value = "'2.319"
$ws.Range("E25").Value = "  +2.26%  "

# Row 26
$ws.Range("D26").Value = "'2.655"
$ws.Range("E26").Value = "  +3.34%  "

# Row 27
$ws.Range("D27").Value = "2.141.17"
$ws.Range("E27").Value = "  +2.76%  "

# Row 28
$ws.Range("D28").Value = "'20.96"
$ws.Range("E28").Value = "  -1.67%  "

# Row 29
$ws.Range("D29").Value = "'158.04"
$ws.Range("E29").Value = "  -0.23%  "

# Row 30
$ws.Range("D30").Value = "'127.03"
$ws.Range("E30").Value = "  -0.70%  "

# Row 31
$ws.Range("D31").Value = "'1.094"
$ws.Range("E31").Value = "  +2.96%  "

# Row 32
$ws.Range("D32").Value = "'0.1066"
$ws.Range("E32").Value = "  +0.17%  "

# Row 33
$ws.Range("D33").Value = "'5.642"
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("D34").Value = "'3.627"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35
$ws.Range("D35").Value = "'9.661"
$ws.Range("E35").Value = "  +1.45%  "

# Row 36
$ws.Range("D36").Value = "'0.06671"
$ws.Range("E36").Value = "  -2.20%  "

# Row 37
$ws.Range("D37").Value = "'0.02421"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").Value = "'1.239"
$ws.Range("E38").Value = "  +1.55%  "

# Row 39
$ws.Range("D39").Value = "'0.2194"
$ws.Range("E39").Value = "  +0.29%  "

# Row 40
$ws.Range("D40").Value = "'1.299"
$ws.Range("E40").Value = "  +10.82%  "

# Row 41
$ws.Range("D41").Value = "'0.6390"
$ws.Range("E41").Value = "  +1.26%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.48"
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'5.002"
$ws.Range("E43").Value = "  -0.35%  "

# Row 44
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("D45").Value = "'13.46"
$ws.Range("E45").Value = "  -1.03%  "

# Row 46
$ws.Range("D46").Value = "'0.6034"
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("D47").Value = "'3.748"
$ws.Range("E47").Value = "  +2.25%  "

# Row 48
$ws.Range("D48").Value = "'1.288"
$ws.Range("E48").Value = "  +1.44%  "

# Row 49
$ws.Range("D49").Value = "'2.042"
$ws.Range("E49").Value = "  +2.43%  "

# Row 50
$ws.Range("D50").Value = "'123.53"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("E51").Value = "  -1.11%  "

